$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.506.42'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.793.21'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  -0.97%  '
$ws.Range("D5").Value = "'339.68"
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("D6").Value = "'0.9983"
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("D7").Value = "'0.3922"
$ws.Range("E7").Value = '  +3.76%  '
$ws.Range("D8").Value = "'0.3464"
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = "'48.22"
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = "'1.193"
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").Value = "'0.07491"
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").Value = "'0.9985"
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = "'21.90"
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").Value = "'6.508"
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = '1.795.83'
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = "'7.151"
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").Value = "'0.00001099"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = "'0.06680"
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = "'84.88"
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").Value = "'0.9975"
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").Value = "'17.69"
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").Value = "'6.550"
$ws.Range("E22").Value = '  +1.78%  '
$ws.Range("D23").Value = '27.513.03'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = "'12.44"
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").Value = "'2.396"
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = "'21.24"
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = "'2.500"
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = "'1.461"
$ws.Range("E28").Value = '  +1.26%  '
$ws.Range("D29").Value = "'155.90"
$ws.Range("E29").Value = '  +4.05%  '
$ws.Range("D30").Value = '1.997.50'
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("D31").Value = "'135.32"
$ws.Range("E31").Value = '  +0.79%  '
$ws.Range("D32").Value = "'4.050"
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").Value = "'6.034"
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").Value = "'0.08773"
$ws.Range("E34").Value = '  +1.54%  '
$ws.Range("D35").Value = "'13.07"
$ws.Range("E35").Value = '  -1.86%  '
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = "'1.617"
$ws.Range("E36").Value = '  -3.43%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = "'5.449"
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = "'0.02411"
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("D39").Value = "'0.06471"
$ws.Range("E39").Value = '  +1.87%  '
$ws.Range("D40").Value = "'0.6814"
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = "'0.2211"
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("D42").Value = "'1.257"
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").Value = "'8.360"
$ws.Range("E43").Value = '  -5.20%  '
$ws.Range("D44").Value = "'14.52"
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = "'0.9966"
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.6398"
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").Value = "'3.875"
$ws.Range("E47").Value = '  +1.36%  '
$ws.Range("D48").Value = "'2.137"
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").Value = "'131.98"
$ws.Range("E49").Value = '  +0.76%  '
$ws.Range("D50").Value = "'0.07196"
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").Value = "'79.77"
$ws.Range("E51").Value = '  +0.43%  '
